$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '34.662.52'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.44%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.790.58'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.50%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '223.51'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.87%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.553'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.11%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.22'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.24%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.27%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0706'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +6.14%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0934'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.67%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.045.82'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.67%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.01'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -4.15%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.784.30'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.94%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.633'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.26%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '34.656.09'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.36%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.30'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.33%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.08'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.49%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '254.02'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.50%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0805'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +7.51%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.999'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.09%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.57'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.24%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.18'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.08%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.05%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '160.20'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.90%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '16.33'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.08%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.11'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.58%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.114'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.96%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.10%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0527'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.14%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.78'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -4.54%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.20'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.86%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.61'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.75%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.87'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.84%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.435.13'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -4.58%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0191'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.79%  '

$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.05'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.16%  '

$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.637'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.42%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '84.73'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.71%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.78'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.89%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.923'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.23%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.32'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.17%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.08'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.94%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.96'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.55%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.87%  '

$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.945.16'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.60%  '

$ws.Range("B47").Value = 'Kaspa'
$ws.Range("C47").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0487'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -6.31%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '105.59'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +7.15%  '

$ws.Range("B49").Value = 'PaxDollar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.999'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.09%  '

$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '11.93'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.50%  '

$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₆0125'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +8.02%  '
